$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two new columns ("gas" after location_(lake), and "measurement" after
# surface_class); downstream columns shift right. It also grows from 6 data rows to 7,
# because each original measurement row (which combined CH4 and CO2 readings) is split
# into one row per gas, pulling the gas type into its own column.

# Insert the two new columns, shifting existing columns to the right.
$ws.Columns("D").Insert()   # new column D = gas
$ws.Columns("I").Insert()   # new column I = measurement

# Clear old data rows; they will be rewritten below, split one-gas-per-row.
$ws.Range("A2:BB7").ClearContents()

# ---- Header row ----
$ws.Range("A1").Value = "program_run?"
$ws.Range("B1").Value = "date_(yyyy-mm-dd)"
$ws.Range("C1").Value = "location_(lake)"
$ws.Range("D1").Value = "gas"
$ws.Range("E1").Value = "start_time_(hh:mm:ss)"
$ws.Range("F1").Value = "stop_time_(hh:mm:ss)"
$ws.Range("G1").Value = "surface_type"
$ws.Range("H1").Value = "surface_class"
$ws.Range("I1").Value = "measurement"
$ws.Range("J1").Value = "measurement_device"
$ws.Range("K1").Value = "collar"
$ws.Range("L1").Value = "collar_height(cm)"
$ws.Range("M1").Value = "submerged_depth(cm)"
$ws.Range("N1").Value = "exposed_height(cm)"
$ws.Range("O1").Value = "Sample ID"
$ws.Range("P1").Value = "Latitude (dd) N +- 3 m"
$ws.Range("Q1").Value = "Longitude (dd) W +- 3m"
$ws.Range("R1").Value = "Waypoint Reference"
$ws.Range("S1").Value = "Position in transect"
$ws.Range("T1").Value = "CH4 flux μmol m^-2 s^-1"
$ws.Range("U1").Value = "CH4 flux ± uncertainty"
$ws.Range("V1").Value = "R_value_used"
$ws.Range("W1").Value = "Use Data? (See Notes)"
$ws.Range("X1").Value = "CO2 Flux μmol m^-2 s^-1"
$ws.Range("Y1").Value = "CO2 flux ± uncertainty"
$ws.Range("Z1").Value = "R_value"
$ws.Range("AA1").Value = "air_Pa"
$ws.Range("AB1").Value = "air_p_mean_Pa"
$ws.Range("AC1").Value = "est_distance_LS (m)"
$ws.Range("AD1").Value = "Submersion Depth (cm)"
$ws.Range("AE1").Value = "Water Temperature (C°)"
$ws.Range("AF1").Value = "notes"
$ws.Range("AG1").Value = "Soil Moisture (VWC %)"
$ws.Range("AH1").Value = "Water pH"
$ws.Range("AI1").Value = "North Soil Temp (C°) 1 cm"
$ws.Range("AJ1").Value = "East Soil Temp (C°) 1 cm"
$ws.Range("AK1").Value = "West Soil Temp (C°) 1 cm"
$ws.Range("AL1").Value = "Average 1 cm Soil Temp (C°)"
$ws.Range("AM1").Value = "North Soil Temp (C°) 5 cm"
$ws.Range("AN1").Value = "East Soil Temp (C°) 5 cm"
$ws.Range("AO1").Value = "West Soil Temp (C°) 5 cm"
$ws.Range("AP1").Value = "Average 5 cm Soil Temp (C°)"
$ws.Range("AQ1").Value = "North Soil Temp (C°) 10 cm"
$ws.Range("AR1").Value = "East Soil Temp (C°) 10 cm"
$ws.Range("AS1").Value = "West Soil Temp (C°) 10 cm"
$ws.Range("AT1").Value = "Average 10 cm Soil Temp (C°)"
$ws.Range("AU1").Value = "North pH"
$ws.Range("AV1").Value = "East pH"
$ws.Range("AW1").Value = "West pH"
$ws.Range("AX1").Value = "Average pH"
$ws.Range("AY1").Value = "Disturbed Surface? (yes) or (no)"
$ws.Range("AZ1").Value = "Overflow Depth (cm)"
$ws.Range("BA1").Value = "Kestrel Data Downloaded?"
$ws.Range("BB1").Value = "GPS data downloaded?"

# ---- Data rows 2:7 ----
# Row 2
$ws.Range("A2").Value = "y"
$ws.Range("B2").Value = "2020-07-30"
$ws.Range("C2").Value = "vault-lake"
$ws.Range("D2").Value = "CH4"
$ws.Range("E2").Value = "17:12:20"
$ws.Range("F2").Value = "17:15:30"
$ws.Range("J2").Value = "bucket"
$ws.Range("K2").Value = "y"
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 34.5
$ws.Range("W2").Value = "rejected"
# Row 3
$ws.Range("A3").Value = "y"
$ws.Range("AA3").Value = 0.9765988463627656
$ws.Range("B3").Value = "2020-07-30"
$ws.Range("C3").Value = "vault-lake"
$ws.Range("D3").Value = "CO2"
$ws.Range("E3").Value = "17:12:20"
$ws.Range("F3").Value = "17:15:30"
$ws.Range("J3").Value = "bucket"
$ws.Range("K3").Value = "y"
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 34.5
$ws.Range("O3").Value = "2020_07_30_17h12m20s_vault-lake_bucket_CO2"
$ws.Range("V3").Value = 0.9817967899243414
$ws.Range("X3").Value = 0.01008573325607265
# Row 4
$ws.Range("A4").Value = "y"
$ws.Range("AA4").Value = 0.9765988463627656
$ws.Range("B4").Value = "2020-07-30"
$ws.Range("C4").Value = "vault-lake"
$ws.Range("D4").Value = "CH4"
$ws.Range("E4").Value = "18:12:20"
$ws.Range("F4").Value = "18:15:30"
$ws.Range("J4").Value = "chamber"
$ws.Range("K4").Value = "y"
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 34.5
$ws.Range("O4").Value = "2020_07_30_18h12m20s_vault-lake_chamber_CH4"
$ws.Range("T4").Value = 171.281040680408
$ws.Range("U4").Value = 0.01693298294852564
$ws.Range("V4").Value = 0.9564391145569133
# Row 5
$ws.Range("A5").Value = "p"
$ws.Range("B5").Value = "2020-07-30"
$ws.Range("C5").Value = "vault-lake"
$ws.Range("D5").Value = "CO2"
$ws.Range("E5").Value = "19:12:20"
$ws.Range("F5").Value = "19:15:30"
$ws.Range("J5").Value = "chamber"
$ws.Range("K5").Value = "n"
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 34.5
# Row 6
$ws.Range("A6").Value = "y"
$ws.Range("B6").Value = "2020-07-30"
$ws.Range("C6").Value = "vault-lake"
$ws.Range("D6").Value = "CH4"
$ws.Range("E6").Value = "20:12:20"
$ws.Range("F6").Value = "20:15:30"
$ws.Range("J6").Value = "chamber"
$ws.Range("K6").Value = "n"
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 34.5
# Row 7
$ws.Range("B7").Value = "2020-07-30"
$ws.Range("C7").Value = "vault-lake"
$ws.Range("D7").Value = "CO2"
$ws.Range("E7").Value = "21:12:20"
$ws.Range("F7").Value = "21:15:30"
$ws.Range("J7").Value = "chamber"
$ws.Range("K7").Value = "n"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 34.5
